# Apply crypto price/volume updates per commit "Updated cryptos list on Sat Sep 21 08:46:18 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.060.42'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.553.09'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '584.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.54'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.36'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('D14').Value = '3.006.71'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '62.935.53'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '2.555.47'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.37'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '336.50'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.78'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.62'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('B27').Value = 'SuiNetwork'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.49'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.40'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.43'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.09%  '
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('D31').Value = '0.0₃0816'
$ws.Range('E31').Value = '  -2.23%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '177.70'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '416.28'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '19.17'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.401'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.48%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.36'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '151.15'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.08%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.77'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '20.92'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0542'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0239'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.35'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('E50').Value = '  -5.36%  '
$ws.Range('E51').Value = '  -0.03%  '
